$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new NPS response rows (391-415 => data rows 390-414)
$ws.Range("A391").Value = 390
$ws.Range("B391").Value = "Pós OS"
$ws.Range("C391").Value = 8005230093
$ws.Range("D391").Value = 46062.708587962959
$ws.Range("E391").Value = "FRQ_ECO_SP_ZONA_SUL_03"
$ws.Range("F391").Value = "Neutros"
$ws.Range("G391").Value = "Boa. As vezes não dá para conciliar disponibilidade de visita"
$ws.Range("H391").Value = "Capacidade"
$ws.Range("I391").Value = "Agenda distante"

$ws.Range("A392").Value = 391
$ws.Range("B392").Value = "Instalação"
$ws.Range("C392").Value = 8005284763
$ws.Range("D392").Value = 46062.929074074083
$ws.Range("E392").Value = "FRQ_ECO_MG_B HORIZONTE"
$ws.Range("F392").Value = "Detratores"
$ws.Range("G392").Value = "A água está com um gosto (horrível) forte de “plástico com bicarbonato”. Está praticamente impossível beber. Já tirei litros e litros e não sai o gosto."
$ws.Range("H392").Value = "Qualidade do Produto"
$ws.Range("I392").Value = "Desempenho da filtragem"

$ws.Range("A393").Value = 392
$ws.Range("B393").Value = "Instalação"
$ws.Range("C393").Value = 8005290922
$ws.Range("D393").Value = 46063.387881944444
$ws.Range("E393").Value = "FRQ_ECO_SP_CAMPINAS_2"
$ws.Range("F393").Value = "Neutros"
$ws.Range("G393").Value = "Purificador foi entregue com vazamento na conexão, após o registro. Ao abrir a solicitação com o time de vendas, não suportam, ao ligar na central após uma ura grande, consegui falar, contudo, não temos site, app ou forma para acompanhar o serviço, sendo a recomendação trazida, até que o ajuste possa ser feito (dia 11/2, sendo que foi instalado dia 06/2), foi de estrangular a mangueira, algo que aumentaria a pressão e não diminuiria a mesma. Sinceramente, a experiência de compra foi muito bom, de instalação razoável, pois ocorreu atraso e não somos informados, e o pós vendas deixando a desejar."
$ws.Range("H393").Value = "Campo"
$ws.Range("I393").Value = "Qualidade da instalação"

$ws.Range("A394").Value = 393
$ws.Range("B394").Value = "Pós OS"
$ws.Range("C394").Value = 8005286290
$ws.Range("D394").Value = 46063.43414351852
$ws.Range("E394").Value = "FRQ_ECO_DF_BRASILIA_2"
$ws.Range("F394").Value = "Neutros"
$ws.Range("G394").Value = "O atendimento é excelente, extremamente profissional. Entretanto, o que tem incomodado é que horário marcado é sistematicamente descumprido, inclusive invertendo o turno solicitado. Reforço que o técnico presta uma excelente serviço."
$ws.Range("H394").Value = "Campo"
$ws.Range("I394").Value = "Fora do período agendado"

$ws.Range("A395").Value = 394
$ws.Range("B395").Value = "Instalação"
$ws.Range("C395").Value = 8005292540
$ws.Range("D395").Value = 46063.443912037037
$ws.Range("E395").Value = "FRQ_ECO_MG_B HORIZONTE"
$ws.Range("F395").Value = "Neutros"
$ws.Range("G395").Value = "A instalação foi super tranquila!"
$ws.Range("H395").Value = "Outros"
$ws.Range("I395").Value = "Satisfação geral"

$ws.Range("A396").Value = 395
$ws.Range("B396").Value = "Pós OS"
$ws.Range("C396").Value = 8005285322
$ws.Range("D396").Value = 46063.468449074076
$ws.Range("E396").Value = "FRQ_ECO_DF_BRASILIA_2"
$ws.Range("F396").Value = "Neutros"
$ws.Range("G396").Value = "gostamos do purificador mas ele é caro e estamos com dificuldade na cobrança"
$ws.Range("H396").Value = "Outros"
$ws.Range("I396").Value = "Preço elevado"

$ws.Range("A397").Value = 396
$ws.Range("B397").Value = "Pós OS"
$ws.Range("C397").Value = 8005296436
$ws.Range("D397").Value = 46063.501562500001
$ws.Range("E397").Value = "FRQ_ECO_SP_S B CAMPO"
$ws.Range("F397").Value = "Detratores"
$ws.Range("G397").Value = "Ainda não resolveu o meu problema."
$ws.Range("H397").Value = "Campo"
$ws.Range("I397").Value = "Reincidência"

$ws.Range("A398").Value = 397
$ws.Range("B398").Value = "Pós OS"
$ws.Range("C398").Value = 8005293020
$ws.Range("D398").Value = 46063.521145833343
$ws.Range("E398").Value = "FRQ_ECO_SP_OSASCO"
$ws.Range("F398").Value = "Detratores"
$ws.Range("G398").Value = "Foram feitos 3 ou 4 agendamentos anteriores onde ninguém apareceu. Nem avisou que não iriam. Precisei agendar novamente todas as outras vezes. Até conseguir esse último tecnico que foi essa semana. Alem disso o tecnico foi no periodo da manhã e estava agendado de tarde. Por sorte tinha gente em casa, senão seria mais um “cano”"
$ws.Range("H398").Value = "Campo"
$ws.Range("I398").Value = "Técnico não cumpriu a agenda"

$ws.Range("A399").Value = 398
$ws.Range("B399").Value = "Pós OS"
$ws.Range("C399").Value = 8005299854
$ws.Range("D399").Value = 46063.522476851853
$ws.Range("E399").Value = "FRQ_ECO_RJ_OESTE"
$ws.Range("F399").Value = "Neutros"
$ws.Range("G399").Value = "O vazamento na parte de baixo esta sendo recorrente.Mais uma vez terá que trocar o aparelho.Ainda não foi feita a troca e continua vazando."
$ws.Range("H399").Value = "Campo"
$ws.Range("I399").Value = "Reincidência"

$ws.Range("A400").Value = 399
$ws.Range("B400").Value = "Pós OS"
$ws.Range("C400").Value = 8005297911
$ws.Range("D400").Value = 46063.522777777784
$ws.Range("E400").Value = "FRQ_ECO_RJ_OESTE"
$ws.Range("F400").Value = "Detratores"
$ws.Range("G400").Value = "O técnico esteve aqui e disse que em 05 dias viria alguém para trocar o purificador. Ninguém veio. Péssimo atendimento."
$ws.Range("H400").Value = "Campo"
$ws.Range("I400").Value = "Reincidência"

$ws.Range("A401").Value = 400
$ws.Range("B401").Value = "Pós OS"
$ws.Range("C401").Value = 8005263938
$ws.Range("D401").Value = 46063.591886574082
$ws.Range("E401").Value = "FRQ_ECO_SP_OSASCO"
$ws.Range("F401").Value = "Neutros"
$ws.Range("G401").Value = "acho muito bom mas acho que vou cancelar a mensalidade esta alta pra mim se houver possibilidade de acordo podemos conversar"
$ws.Range("H401").Value = "Outros"
$ws.Range("I401").Value = "Preço elevado"

$ws.Range("A402").Value = 401
$ws.Range("B402").Value = "Pós OS"
$ws.Range("C402").Value = 8005286964
$ws.Range("D402").Value = 46063.629861111112
$ws.Range("E402").Value = "AT_ECO_SP_BAURU"
$ws.Range("F402").Value = "Neutros"
$ws.Range("G402").Value = "Boa tarde. O técnico foi atenciosao e veio dentro do horario marcado.Preciso que o técnico retorne ao meu apto, o purificador está fazendo barulho quando acionado.Outra coisa que me chamou atenção,  foi que ele fez a troca do purificador mas não trocou o elemento filtrante externo.Moro em um prédio com tubulação em ferro e só o interno não é suficiente.Aguardo manifestação de vcs."
$ws.Range("H402").Value = "Campo"
$ws.Range("I402").Value = "Qualidade da manutenção"

$ws.Range("A403").Value = 402
$ws.Range("B403").Value = "Pós OS"
$ws.Range("C403").Value = 8005292458
$ws.Range("D403").Value = 46063.795659722222
$ws.Range("E403").Value = "FRQ_ECO_SP_OSASCO"
$ws.Range("F403").Value = "Neutros"
$ws.Range("G403").Value = "Nao tenho problema com a qualidade da agua. Os aparelhos atuais dao problemas mecânicos mais rapido que o antigo."
$ws.Range("H403").Value = "Qualidade do Produto"
$ws.Range("I403").Value = "Funcionamento geral"

$ws.Range("A404").Value = 403
$ws.Range("B404").Value = "Pós OS"
$ws.Range("C404").Value = 8005284928
$ws.Range("D404").Value = 46063.810046296298
$ws.Range("E404").Value = "FRQ_ECO_RJ_OESTE"
$ws.Range("F404").Value = "Detratores"
$ws.Range("G404").Value = "Acho que o valor do serviço é muito alto"
$ws.Range("H404").Value = "Outros"
$ws.Range("I404").Value = "Preço elevado"

$ws.Range("A405").Value = 404
$ws.Range("B405").Value = "Pós OS"
$ws.Range("C405").Value = 8005262658
$ws.Range("D405").Value = 46064.367974537039
$ws.Range("E405").Value = "FRQ_ECO_SP_S B CAMPO"
$ws.Range("F405").Value = "Neutros"
$ws.Range("G405").Value = "Demorou um mês entre minha ligação e o agendamento da visita, é a primeira vez que demora este período longo."
$ws.Range("H405").Value = "Capacidade"
$ws.Range("I405").Value = "Agenda distante"

$ws.Range("A406").Value = 405
$ws.Range("B406").Value = "Pós OS"
$ws.Range("C406").Value = 8005286061
$ws.Range("D406").Value = 46064.431423611109
$ws.Range("E406").Value = "FRQ_ECO_RJ_OESTE"
$ws.Range("F406").Value = "Neutros"
$ws.Range("G406").Value = "O Técnico habitual está de férias e o seu substituto deixou a desejar.Em 5 minutos trocou o filtro e não testou quimicamente a água,  nem purgou o ar do aparelho.Ao chegar, não se apresentou; não reparou na abertura da válvula de água  e deixou o material utilizado espalhados sobre o balcão. Creio que, por ser novato e por cobrir a ausência daquele que está de férias, não julgou necessário valorizar o atendimento. Tomara que a experiência lhe traga mais comprometimento com a imagem da companhia."
$ws.Range("H406").Value = "Campo"
$ws.Range("I406").Value = "Competência do técnico"

$ws.Range("A407").Value = 406
$ws.Range("B407").Value = "Pós OS"
$ws.Range("C407").Value = 8005298244
$ws.Range("D407").Value = 46064.435601851852
$ws.Range("E407").Value = "FRQ_ECO_RJ_OESTE"
$ws.Range("F407").Value = "Detratores"
$ws.Range("G407").Value = "Como vou avaliar uma coisa que não teve.Atendimento pessímo, não indico."
$ws.Range("H407").Value = "Campo"
$ws.Range("I407").Value = "Técnico não cumpriu a agenda"

$ws.Range("A408").Value = 407
$ws.Range("B408").Value = "Instalação"
$ws.Range("C408").Value = 8005297440
$ws.Range("D408").Value = 46064.474849537037
$ws.Range("E408").Value = "FRQ_ECO_RJ_ATLANTICA"
$ws.Range("F408").Value = "Detratores"
$ws.Range("G408").Value = "Eu estava em casa o rapaz nao veio e depois quando veio disse que precisa de um pressurizador em vez da Brastemp já deixar com o tecnico precisa aguardar mais de 10 dias. Sinceramente nao tem 1 semana e ja me arrependi amargamente"
$ws.Range("H408").Value = "Campo"
$ws.Range("I408").Value = "Técnico não cumpriu a agenda"

$ws.Range("A409").Value = 408
$ws.Range("B409").Value = "Pós OS"
$ws.Range("C409").Value = 8005269484
$ws.Range("D409").Value = 46064.502465277779
$ws.Range("E409").Value = "FRQ_ECO_SP_S B CAMPO"
$ws.Range("F409").Value = "Neutros"
$ws.Range("G409").Value = "Gosto muito, só acho o valor um pouco alto e demora para agendar visita técnica"
$ws.Range("H409").Value = "Capacidade"
$ws.Range("I409").Value = "Agenda distante"

$ws.Range("A410").Value = 409
$ws.Range("B410").Value = "Pós OS"
$ws.Range("C410").Value = 8005290057
$ws.Range("D410").Value = 46064.505196759259
$ws.Range("E410").Value = "FRQ_ECO_RJ_OESTE"
$ws.Range("F410").Value = "Detratores"
$ws.Range("G410").Value = "Pela pontualidade, foi marcado na parte da manhã , isto é, 8 às 13h. Fiquei esperando , amanhã toda. Abrir mão de outras coisas, para ficar a disposição. Só chegou após 13:20h. Estou insatisfeito com atendimento."
$ws.Range("H410").Value = "Campo"
$ws.Range("I410").Value = "Fora do período agendado"

$ws.Range("A411").Value = 410
$ws.Range("B411").Value = "Pós OS"
$ws.Range("C411").Value = 8005285756
$ws.Range("D411").Value = 46064.515173611107
$ws.Range("E411").Value = "FRQ_ECO_RJ_ATLANTICA"
$ws.Range("F411").Value = "Neutros"
$ws.Range("G411").Value = "Purificador fornece água bem filtrada,límpida,sem gosto de cloro.O problema é o atendimento virtual.Quando liga para marcar pergunta se é o responsável,resposta SIM.Pergunta novamente,várias vezes,e nada.Desliga.Quando novo contato,oferece data vencida.Até conseguir agendar,toma muito tempo."
$ws.Range("H411").Value = "Atendimento "
$ws.Range("I411").Value = "Erro de comunicação"

$ws.Range("A412").Value = 411
$ws.Range("B412").Value = "Pós OS"
$ws.Range("C412").Value = 8005264422
$ws.Range("D412").Value = 46064.535173611112
$ws.Range("E412").Value = "FRQ_ECO_SP_ZONA_SUL_03"
$ws.Range("F412").Value = "Neutros"
$ws.Range("G412").Value = "Bom"
$ws.Range("H412").Value = "Outros"
$ws.Range("I412").Value = "Satisfação geral"

$ws.Range("A413").Value = 412
$ws.Range("B413").Value = "Pós OS"
$ws.Range("C413").Value = 8005296356
$ws.Range("D413").Value = 46064.545416666668
$ws.Range("E413").Value = "FRQ_ECO_SP_SANTOS"
$ws.Range("F413").Value = "Neutros"
$ws.Range("G413").Value = "Aluguei o purificador recentemente (2 meses) e o aparelho apresentou defeito, segundo o técnico, falta de gás. O atendimento do técnico foi muito bom porém o agendamento para a visita e agendamento para substituição do  purificador é demorada. Quando reinstalarem pedirei o abatimento no boleto dos duas que fiquei sem o purificador pois não quero pagar por algo que contratei e a Brastemp não está fornecendo."
$ws.Range("H413").Value = "Qualidade do Produto"
$ws.Range("I413").Value = "Funcionamento geral"

$ws.Range("A414").Value = 413
$ws.Range("B414").Value = "Pós OS"
$ws.Range("C414").Value = 8005273187
$ws.Range("D414").Value = 46064.558206018519
$ws.Range("E414").Value = "FRQ_ECO_SP_OSASCO"
$ws.Range("F414").Value = "Neutros"
$ws.Range("G414").Value = "O técnico não sabia o número da ordem de serviço e nem a senha, o que é muito importante por questões de segurança."
$ws.Range("H414").Value = "Campo"
$ws.Range("I414").Value = "Competência do técnico"

$ws.Range("A415").Value = 414
$ws.Range("B415").Value = "Pós OS"
$ws.Range("C415").Value = 8005263760
$ws.Range("D415").Value = 46064.574062500003
$ws.Range("E415").Value = "FRQ_ECO_SP_OSASCO"
$ws.Range("F415").Value = "Neutros"
$ws.Range("G415").Value = "A água é muito boa, porém em caso de falta de energia elétrica, morremos de sede…."
$ws.Range("H415").Value = "Qualidade do Produto"
$ws.Range("I415").Value = "Funcionamento geral"

# Update the active sheet view / selection to match the saved workbook state
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 381
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E406").Select() | Out-Null
